$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18's phone number (A18) was entered as text; normalize it to a number,
# same as every other "phone" cell in column A.
$ws.Range("A18").Value = 79174445

# Append the new payment row (19) for phone 79174445, paid by Cash.
# Column A keeps the original text formatting this record was logged with,
# so force Text format before writing it (otherwise Excel auto-detects the
# numeric-looking string and stores it as a number).
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "79174445"
$ws.Range("A19").Style = "Normal"

$ws.Range("C19").Value = "Cash"
$ws.Range("D19").Value = "2025-08-18T08:57:38"
$ws.Range("E19").Value = 30
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 30
